$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing header cell (H1) into the two new
# header cells so they match the style used by the other headers (bold,
# centered, bordered).
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

# Set header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new data columns I and J for rows 2-11
$data = @(
    @(6, 6),
    @(5, 6),
    @(4, 6),
    @(5, 8),
    @(7, 8),
    @(7, 9),
    @(7, 9),
    @(4, 5),
    @(5, 5),
    @(4, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
